$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Worksheet, $CellRef, $TextValue)
    $range = $Worksheet.Range($CellRef)
    $range.NumberFormat = "@"
    $range.Value = $TextValue
    $range.Style = "Normal"
}

Set-TextValue $ws "D2" "42.864.04"
Set-TextValue $ws "E2" "  -2.28%  "
Set-TextValue $ws "D3" "2.556.15"
Set-TextValue $ws "E3" "  -1.89%  "
Set-TextValue $ws "D4" "0.999"
Set-TextValue $ws "E4" "  +0.05%  "
Set-TextValue $ws "D5" "301.90"
Set-TextValue $ws "E5" "  +0.04%  "
Set-TextValue $ws "D6" "94.68"
Set-TextValue $ws "E6" "  -2.38%  "
Set-TextValue $ws "E7" "  -0.88%  "
Set-TextValue $ws "E8" "  +0.03%  "
Set-TextValue $ws "D9" "0.544"
Set-TextValue $ws "E9" "  -2.98%  "
Set-TextValue $ws "D10" "36.30"
Set-TextValue $ws "E10" "  -2.75%  "
Set-TextValue $ws "D11" "0.0810"
Set-TextValue $ws "E11" "  -0.80%  "
Set-TextValue $ws "D12" "7.77"
Set-TextValue $ws "E12" "  -1.39%  "
Set-TextValue $ws "E13" "  +6.19%  "
Set-TextValue $ws "D14" "2.544.55"
Set-TextValue $ws "E14" "  -2.21%  "
Set-TextValue $ws "E15" "  -1.69%  "
Set-TextValue $ws "D16" "14.21"
Set-TextValue $ws "E16" "  -1.73%  "
Set-TextValue $ws "D17" "42.919.81"
Set-TextValue $ws "E17" "  -2.05%  "
Set-TextValue $ws "D18" "0.0₃0988"
Set-TextValue $ws "E18" "  +0.78%  "
Set-TextValue $ws "E19" "  +1.45%  "
Set-TextValue $ws "D20" "6.57"
Set-TextValue $ws "E20" "  -1.74%  "
Set-TextValue $ws "D21" "71.56"
Set-TextValue $ws "E21" "  -2.41%  "
Set-TextValue $ws "D22" "252.41"
Set-TextValue $ws "E22" "  -5.58%  "
Set-TextValue $ws "D23" "2.93"
Set-TextValue $ws "E23" "  -0.60%  "
Set-TextValue $ws "E24" "  -5.16%  "
Set-TextValue $ws "D25" "28.56"
Set-TextValue $ws "E25" "  -4.01%  "
Set-TextValue $ws "D26" "0.999"
Set-TextValue $ws "E26" "  -0.21%  "
Set-TextValue $ws "D27" "10.23"
Set-TextValue $ws "E27" "  -0.85%  "
Set-TextValue $ws "D28" "36.95"
Set-TextValue $ws "E28" "  -2.78%  "
Set-TextValue $ws "D29" "2.12"
Set-TextValue $ws "E29" "  -1.32%  "
Set-TextValue $ws "E30" "  -0.87%  "
Set-TextValue $ws "D31" "153.60"
Set-TextValue $ws "E31" "  +1.05%  "
Set-TextValue $ws "E32" "  -1.48%  "
Set-TextValue $ws "E33" "  -6.81%  "
Set-TextValue $ws "D34" "2.14"
Set-TextValue $ws "E34" "  -5.91%  "
Set-TextValue $ws "D35" "0.0798"
Set-TextValue $ws "E35" "  -2.28%  "
Set-TextValue $ws "D36" "18.10"
Set-TextValue $ws "E36" "  +6.81%  "
Set-TextValue $ws "D37" "0.113"
Set-TextValue $ws "E37" "  -3.64%  "
Set-TextValue $ws "D38" "0.119"
Set-TextValue $ws "E38" "  -1.42%  "
Set-TextValue $ws "D39" "23.45"
Set-TextValue $ws "E39" "  -4.23%  "
Set-TextValue $ws "D40" "2.13"
Set-TextValue $ws "E40" "  +32.57%  "
Set-TextValue $ws "D41" "3.41"
Set-TextValue $ws "E41" "  -3.90%  "
Set-TextValue $ws "B42" "VeChain"
Set-TextValue $ws "C42" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws "D42" "0.0310"
Set-TextValue $ws "E42" "  -2.04%  "
Set-TextValue $ws "B43" "RenderToken"
Set-TextValue $ws "C43" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws "D43" "3.87"
Set-TextValue $ws "E43" "  -0.30%  "
Set-TextValue $ws "D44" "2.086.81"
Set-TextValue $ws "E44" "  +0.55%  "
Set-TextValue $ws "E45" "  +0.36%  "
Set-TextValue $ws "D46" "9.26"
Set-TextValue $ws "E46" "  +0.41%  "
Set-TextValue $ws "D47" "84.76"
Set-TextValue $ws "E47" "  -4.51%  "
Set-TextValue $ws "B48" "Aave"
Set-TextValue $ws "C48" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws "D48" "106.55"
Set-TextValue $ws "E48" "  -0.13%  "
Set-TextValue $ws "B49" "RocketPoolETH"
Set-TextValue $ws "C49" "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue $ws "D49" "2.805.38"
Set-TextValue $ws "E49" "  -1.40%  "
Set-TextValue $ws "B50" "ordi"
Set-TextValue $ws "C50" "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
Set-TextValue $ws "D50" "75.08"
Set-TextValue $ws "E50" "  +7.68%  "
Set-TextValue $ws "B51" "Algorand"
Set-TextValue $ws "C51" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws "D51" "0.191"
Set-TextValue $ws "E51" "  -0.63%  "
